# projeto02 - feat e check
# Add year/month/day metadata rows to the "work" sheet and tidy up the
# trailing helper-style cells that used to sit below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old rows 11/12 only carried leftover cell styling (no data) in H11,
# I11 and G12. Clear that formatting cell-by-cell (instead of over the
# whole future range) so we don't materialize extra empty cells when
# writing the new rows below.
$ws.Range("H11").ClearFormats()
$ws.Range("I11").ClearFormats()
$ws.Range("G12").ClearFormats()

# New metadata rows describing the year/month/day columns added to the
# nyflights table.
$ws.Range("A11").Value = "nyflights"
$ws.Range("B11").Value = "year"
$ws.Range("C11").Value = "int"
$ws.Range("D11").Value = "ano"
$ws.Range("E11").Value = "int"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0.1
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

$ws.Range("A12").Value = "nyflights"
$ws.Range("B12").Value = "month"
$ws.Range("C12").Value = "int"
$ws.Range("D12").Value = "mes"
$ws.Range("E12").Value = "int"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0.1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

$ws.Range("A13").Value = "nyflights"
$ws.Range("B13").Value = "day"
$ws.Range("C13").Value = "int"
$ws.Range("D13").Value = "dia"
$ws.Range("E13").Value = "int"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0.1
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# Row 14 stays blank but keeps the underlined "divider" style that used to
# live a couple of rows further down; row 15 keeps a single styled (but
# empty) cell in column C, mirroring the old C16.
$ws.Range("A14:I14").Font.Underline = $true
$ws.Range("C15").Font.Underline = $true

# G15 / C16 were just leftover stray styled cells below the old table; drop
# them completely now that the table grew to cover rows 11-13.
$ws.Range("G15").ClearFormats()
$ws.Range("G15").ClearContents()
$ws.Range("C16").ClearFormats()
$ws.Range("C16").ClearContents()

$ws.Range("C11:C13").Select() | Out-Null
